# Added PV FK's to schema and updated demo data.
#
# - "person" sheet header row: rename "Email"/"Nick" labels to the new
#   lowercase "email"/"nick" field-key labels (the schema now stores
#   lowercase PV foreign-key names).
# - "person" sheet H column (initials): swap the concatenation order from
#   "first-initial + last-initial" to "last-initial + first-initial".

$wb = $excel.ActiveWorkbook
$person = $wb.Worksheets.Item("person")

# --- Header row relabel -----------------------------------------------
$person.Range("D1").Value = "email"
$person.Range("E1").Value = "nick"

# --- Swap initials formula (LEFT(lastname) & LEFT(firstname)) ---------
$person.Range("H2").Formula = "=CONCATENATE(LEFT(G2,1),LEFT(F2,1))"
$person.Range("H3:H23").Formula = "=CONCATENATE(LEFT(G3,1),LEFT(F3,1))"
